# This script refreshes the LR-pair statistics table (NATMI TPM output)
# for Cthrc1-Ror2 with newly recomputed values, and appends the new
# "Resolving-Mac" sending-cluster rows (11-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: SendingCluster, LigandSymbol, ReceptorSymbol, TargetCluster, then
# the 16 numeric statistic columns (E..T).
$rows = @(
    @(2, "ECs", "Cthrc1", "Ror2", "ECs", 2, 0.6666666666666666, 0.8618196666666668, 2.585459, 0.006773656541421759, 0.006773656541421758, 1, 0.3333333333333333, 0.028814, 0.086442, 0.003707384188741118, 0.003707384188741118, 0.02483247187533334, 0.223492246878, 0.00002511254716162988, 0.00002511254716162987),
    @(3, "ECs", "Cthrc1", "Ror2", "FAPs", 2, 0.6666666666666666, 0.8618196666666668, 2.585459, 0.006773656541421759, 0.006773656541421758, 3, 1, 7.377589, 22.132767, 0.9492453949340737, 0.9492453949340736, 6.358151292783668, 57.22336163505301, 0.00642986227880967, 0.006429862278809667),
    @(4, "ECs", "Cthrc1", "Ror2", "MuSCs", 2, 0.6666666666666666, 0.8618196666666668, 2.585459, 0.006773656541421759, 0.006773656541421758, 3, 1, 0.3656536666666667, 1.096961, 0.04704722087718523, 0.04704722087718523, 0.3151275211221112, 2.836147690099001, 0.0003186817154504601, 0.00031868171545046),
    @(5, "FAPs", "Cthrc1", "Ror2", "ECs", 3, 1, 125.8872733333333, 377.66182, 0.9894380291809874, 0.9894380291809874, 1, 0.3333333333333333, 0.028814, 0.086442, 0.003707384188741118, 0.003707384188741118, 3.627315893826667, 32.64584304444001, 0.003668226905124766, 0.003668226905124766),
    @(6, "FAPs", "Cthrc1", "Ror2", "FAPs", 3, 1, 125.8872733333333, 377.66182, 0.9894380291809874, 0.9894380291809874, 3, 1, 7.377589, 22.132767, 0.9492453949340737, 0.9492453949340736, 928.7445629839934, 8358.701066855941, 0.9392194927726979, 0.9392194927726978),
    @(7, "FAPs", "Cthrc1", "Ror2", "MuSCs", 3, 1, 125.8872733333333, 377.66182, 0.9894380291809874, 0.9894380291809874, 3, 1, 0.3656536666666667, 1.096961, 0.04704722087718523, 0.04704722087718523, 46.03114308100223, 414.2802877290201, 0.04655030950316476, 0.04655030950316476),
    @(8, "MuSCs", "Cthrc1", "Ror2", "ECs", 3, 1, 0.4264636666666666, 1.279391, 0.00335188267003504, 0.00335188267003504, 1, 0.3333333333333333, 0.028814, 0.086442, 0.003707384188741118, 0.003707384188741118, 0.01228812409133333, 0.110593116822, 0.00001242671681340327, 0.00001242671681340327),
    @(9, "MuSCs", "Cthrc1", "Ror2", "FAPs", 3, 1, 0.4264636666666666, 1.279391, 0.00335188267003504, 0.00335188267003504, 3, 1, 7.377589, 22.132767, 0.9492453949340737, 0.9492453949340736, 3.146273656099666, 28.316462904897, 0.003181759188890089, 0.003181759188890088),
    @(10, "MuSCs", "Cthrc1", "Ror2", "MuSCs", 3, 1, 0.4264636666666666, 1.279391, 0.00335188267003504, 0.00335188267003504, 3, 1, 0.3656536666666667, 1.096961, 0.04704722087718523, 0.04704722087718523, 0.1559380034167778, 1.403442030751, 0.0001576967643315479, 0.0001576967643315479),
    @(11, "Resolving-Mac", "Cthrc1", "Ror2", "ECs", 1, 0.3333333333333333, 0.05552766666666667, 0.166583, 0.0004364316075558192, 0.0004364316075558192, 1, 0.3333333333333333, 0.028814, 0.086442, 0.003707384188741118, 0.003707384188741118, 0.001599974187333334, 0.014399767686, 0.000001618019641319313, 0.000001618019641319313),
    @(12, "Resolving-Mac", "Cthrc1", "Ror2", "FAPs", 1, 0.3333333333333333, 0.05552766666666667, 0.166583, 0.0004364316075558192, 0.0004364316075558192, 3, 1, 7.377589, 22.132767, 0.9492453949340737, 0.9492453949340736, 0.4096603027956667, 3.686942725161, 0.0004142806936760363, 0.0004142806936760362),
    @(13, "Resolving-Mac", "Cthrc1", "Ror2", "MuSCs", 1, 0.3333333333333333, 0.05552766666666667, 0.166583, 0.0004364316075558192, 0.0004364316075558192, 3, 1, 0.3656536666666667, 1.096961, 0.04704722087718523, 0.04704722087718523, 0.02030389491811111, 0.182735054263, 0.00002053289423846365, 0.00002053289423846365)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A" + $rowNum).Value = $r[1]
    $ws.Range("B" + $rowNum).Value = $r[2]
    $ws.Range("C" + $rowNum).Value = $r[3]
    $ws.Range("D" + $rowNum).Value = $r[4]
    $ws.Range("E" + $rowNum).Value = $r[5]
    $ws.Range("F" + $rowNum).Value = $r[6]
    $ws.Range("G" + $rowNum).Value = $r[7]
    $ws.Range("H" + $rowNum).Value = $r[8]
    $ws.Range("I" + $rowNum).Value = $r[9]
    $ws.Range("J" + $rowNum).Value = $r[10]
    $ws.Range("K" + $rowNum).Value = $r[11]
    $ws.Range("L" + $rowNum).Value = $r[12]
    $ws.Range("M" + $rowNum).Value = $r[13]
    $ws.Range("N" + $rowNum).Value = $r[14]
    $ws.Range("O" + $rowNum).Value = $r[15]
    $ws.Range("P" + $rowNum).Value = $r[16]
    $ws.Range("Q" + $rowNum).Value = $r[17]
    $ws.Range("R" + $rowNum).Value = $r[18]
    $ws.Range("S" + $rowNum).Value = $r[19]
    $ws.Range("T" + $rowNum).Value = $r[20]
}

Write-Host "Updated $($rows.Count) rows (2..13) with refreshed TPM statistics."